$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where C:F are the "roll-up" (sum of the two detail rows immediately below)
$rows = @(12, 15, 19, 22, 26, 33, 36, 57)

foreach ($r in $rows) {
    $below1 = $r + 1
    $below2 = $r + 2

    $ws.Range("C$r").Formula = "=C$($below1)+C$($below2)"
    $ws.Range("D$r").Formula = "=D$($below1)+D$($below2)"
    $ws.Range("E$r").Formula = "=E$($below1)+E$($below2)"
    $ws.Range("F$r").Formula = "=F$($below1)+F$($below2)"
}

# Row 12's "C" cell did not exist before, so it adopted the column's default
# format (style index 4) when it was created; normalize the rest of that row
# (D12:F12, which previously carried the old styles 9/9/10) to match, the
# same way Excel would if the whole C12:F12 block were formatted together.
$ws.Range("C12").Copy() | Out-Null
$ws.Range("D12:F12").PasteSpecial(-4122) | Out-Null

# For the other roll-up rows, C/D/E already existed with style 9; only F
# (previously style 10, the "with fill" variant) needs to be normalized back
# down to plain style 9 to match its row mates.
foreach ($r in @(15, 19, 22, 26, 33, 36, 57)) {
    $ws.Range("C$r").Copy() | Out-Null
    $ws.Range("F$r").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

$ws.Range("C57:F57").Select()

$wb.Save() | Out-Null
